$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("K1").Value = "next_ex_dividend_dt"
$ws.Range("L1").Value = "dividend_yield"

# Row 2 - AFCG
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 13.399
$ws.Range("H2").Value = 0

# Row 3 - AGNC
$ws.Range("G3").Value = 10.285

# Row 4 - AMZN
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 258
$ws.Range("G4").Value = 130.8301
$ws.Range("H4").Value = 392.4902999999999

# Row 5 - ARLP
$ws.Range("G5").Value = 17.86

# Row 6 - ARR
$ws.Range("G6").Value = 5.34

# Row 7 - BDN
$ws.Range("G7").Value = 4.505

# Row 8 - BRY
$ws.Range("G8").Value = 6.62
$ws.Range("H8").Value = 92.68000000000001

# Row 9 - CIM
$ws.Range("G9").Value = 5.92

# Row 10 - ECC
$ws.Range("C10").Value = "MONITOR"
$ws.Range("G10").Value = 10.355
$ws.Range("H10").Value = 155.325

# Row 11 - EFC
$ws.Range("G11").Value = 13.79
$ws.Range("H11").Value = 82.73999999999999

# Row 12 - EGLE
$ws.Range("C12").Value = "MONITOR"
$ws.Range("G12").Value = 44.81
$ws.Range("H12").Value = 44.81

# Row 13 - FRO
$ws.Range("G13").Value = 13.5

# Row 14 - GGB
$ws.Range("G14").Value = 5.15
$ws.Range("H14").Value = 77.25

# Row 15 - GNK
$ws.Range("G15").Value = 13.3539
$ws.Range("H15").Value = 66.76949999999999

# Row 16 - GOGL
$ws.Range("G16").Value = 7.3
$ws.Range("H16").Value = 65.7

# Row 17 - GOOGL
$ws.Range("C17").Value = "MONITOR"
$ws.Range("G17").Value = 119.87
$ws.Range("H17").Value = 479.48

# Row 18 - ICMB
$ws.Range("G18").Value = 3.65

# Row 19 - KEN
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = 24.376
$ws.Range("F19").Value = 804.408
$ws.Range("G19").Value = 23.84
$ws.Range("H19").Value = 786.72
$ws.Range("I19").Value = 24.376
$ws.Range("J19").Value = 28.03

# Row 20 - LND
$ws.Range("C20").Value = "MONITOR"
$ws.Range("D20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 5.0799
$ws.Range("H20").Value = 0

# Row 21 - MFA
$ws.Range("G21").Value = 11.44

# Row 22 - NLY
$ws.Range("G22").Value = 20.59

# Row 23 - NYMT
$ws.Range("G23").Value = 9.9

# Row 24 - OPI
$ws.Range("G24").Value = 7.89

# Row 25 - ORC
$ws.Range("G25").Value = 10.3675

# Row 26 - OXLC
$ws.Range("C26").Value = "MONITOR"
$ws.Range("D26").Value = 366
$ws.Range("E26").Value = 4.916
$ws.Range("F26").Value = 1799.256
$ws.Range("G26").Value = 4.975
$ws.Range("H26").Value = 1820.85
$ws.Range("I26").Value = 4.916
$ws.Range("J26").Value = 5.65

# Row 27 - PDM
$ws.Range("D27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 7.23
$ws.Range("H27").Value = 0

# Row 28 - PTMN
$ws.Range("G28").Value = 19.64

# Row 29 - SBLK
$ws.Range("G29").Value = 17.125
$ws.Range("H29").Value = 85.625

# Row 30 - SOFI
$ws.Range("G30").Value = 9.037599999999999
$ws.Range("H30").Value = 9.037599999999999

# Row 31 - TRTX
$ws.Range("G31").Value = 7.345

# Row 32 - TSLA
$ws.Range("G32").Value = 257.5918
$ws.Range("H32").Value = 257.5918

# Row 33 - TWO
$ws.Range("G33").Value = 13.605
